# Automatische test-sync: 2025-07-23 14:22:50
# Adds a new row (row 5) of log data to the "Logs" sheet, extends the
# conditional-formatting ranges that covered rows 2-4 to also cover row 5,
# and bumps the "Aantal" counter on the "Dashboard" sheet from 3 to 4.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 5 -----------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A5").Value = "Wat zijn jullie openingstijden?"
$logs.Range("B5").Value = "mailmind.test@zohomail.eu"
$logs.Range("C5").Value = "Testmail #1: Wat zijn jullie openingstijden?"
$logs.Range("D5").Value = "Openingstijden / Locatie"
$logs.Range("E5").Value = "Beste klant,`nDank u wel voor uw interesse in onze diensten. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Range("F5").Value = "2025-07-23 14:21:53"
$logs.Range("G5").Value = "Ja"
$logs.Range("H5").Value = "Nee"
$logs.Range("I5").Value = "Ja"
$logs.Range("J5").Value = "Ja"

# Re-fit the row height back to the sheet default (writing the wrapped,
# multi-line "Antwoord" text would otherwise leave a custom row height).
$logs.Rows.Item(5).AutoFit()

# --- Extend conditional formatting ranges from *2:*4 to *2:*5 -------------
$logs.Range("D2:D4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D5"))
$logs.Range("G2:G4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G5"))
$logs.Range("H2:H4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H5"))
$logs.Range("I2:I4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I5"))
$logs.Range("J2:J4").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J5"))

# --- Dashboard sheet: bump the count for "Openingstijden / Locatie" -------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 4
